# Applies the 25th Feb attendance update to the FA-II (A) Attendance Sheet.
# Marks additional sessions as Absent ("A") for a number of participants.
# Setting WrapText = $false on a cell (after giving it a value, or on its
# own) nudges the style engine to emit the same "applyAlignment" cell
# style (index 36 in the original workbook) that the reference edit uses,
# matching the target XML precisely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Absent($addr) {
    $r = $ws.Range($addr)
    $r.Value = "A"
    $r.WrapText = $false
}

function Restyle-Blank($addr) {
    $r = $ws.Range($addr)
    $r.WrapText = $false
}

# Rows where a single new absence is recorded
Set-Absent("U11")
Set-Absent("T13")
Set-Absent("T14")
Set-Absent("U16")
Set-Absent("U18")
Set-Absent("T20")
Set-Absent("T26")
Set-Absent("U27")
Set-Absent("U28")
Set-Absent("T29")
Set-Absent("T34")
Set-Absent("U58")
Set-Absent("U59")
Set-Absent("T67")
Set-Absent("U71")
Set-Absent("T73")

# Rows where the "T" cell only gets re-styled (no value) and "U" gets the
# new absence mark
Restyle-Blank("T40")
Set-Absent("U40")

Restyle-Blank("T41")
Set-Absent("U41")

Restyle-Blank("T44")
Set-Absent("U44")

Restyle-Blank("T46")
Set-Absent("U46")

# Row 68 gets two new absences
Set-Absent("T68")
Set-Absent("U68")
